$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 67; this shifts the existing rows 67-70 down to 68-71,
# preserving their data (this matches the diff: old row67 data -> new row68,
# old row68 -> new row69, old row69 -> new row70, old row70 -> new row71).
$ws.Rows("67:67").Insert()

# Populate the newly inserted row 67 with this week's new record (same shape as the
# record that used to be at row 67, but with a new "Fecha" date serial of 45106).
$ws.Cells.Item(67, 1).Value = 9
$ws.Cells.Item(67, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(67, 3).Value = "Metropolitana"
$ws.Cells.Item(67, 4).Value = 45106
$ws.Cells.Item(67, 5).Value = 13
$ws.Cells.Item(67, 6).Value = 100112010
$ws.Cells.Item(67, 7).Value = "Achicoria"
$ws.Cells.Item(67, 8).Value = "Sin especificar"
$ws.Cells.Item(67, 9).Value = "Primera"
$ws.Cells.Item(67, 10).Value = 70
$ws.Cells.Item(67, 11).Value = 7000
$ws.Cells.Item(67, 12).Value = 7000
$ws.Cells.Item(67, 13).Value = 7000
$ws.Cells.Item(67, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(67, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(67, 16).Value = 438
$ws.Cells.Item(67, 17).Value = 16
$ws.Cells.Item(67, 18).Value = "Hortaliza"
